# Update countries & provincias Spain
# - Reorders "Togo" earlier in the country list (now appears right after
#   "Mauricio", before "Cabo Verde") and refreshes the covid-19 stats for
#   the affected rows.
# - Updates various case-count figures (Estados Unidos, Alemania, Maldivas,
#   the African-country block rows, Yemen).
# - Refreshes the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 23:05"

# --- simple statistic refreshes (no text change) ---------------------------
# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1543842   # Casos totales
$ws.Cells.Item(4, 3).Value = 16178     # Nuevos casos
$ws.Cells.Item(4, 4).Value = 351936    # Casos activos
$ws.Cells.Item(4, 5).Value = 1100223   # Recuperados
$ws.Cells.Item(4, 7).Value = 705       # Casos criticos
$ws.Cells.Item(4, 8).Value = 91683     # Muertes

# Alemania (row 11)
$ws.Cells.Item(11, 2).Value = 177281
$ws.Cells.Item(11, 3).Value = 630
$ws.Cells.Item(11, 5).Value = 14561
$ws.Cells.Item(11, 7).Value = 71
$ws.Cells.Item(11, 8).Value = 8120

# Maldivas (row 101)
$ws.Cells.Item(101, 2).Value = 1106
$ws.Cells.Item(101, 3).Value = 12
$ws.Cells.Item(101, 5).Value = 1044

# --- Togo moved up in the ranking, shifting the text of rows 140-145 -------
# (numbers for each row are refreshed to the new reported figures as well)

# Row 140 -> now "Togo"
$ws.Cells.Item(140, 1).Value = "Togo"
$ws.Cells.Item(140, 2).Value = 330
$ws.Cells.Item(140, 3).Value = 29
$ws.Cells.Item(140, 4).Value = 106
$ws.Cells.Item(140, 5).Value = 212
$ws.Cells.Item(140, 7).Value = 1
$ws.Cells.Item(140, 8).Value = 12

# Row 141 -> now "Cabo Verde"
$ws.Cells.Item(141, 1).Value = "Cabo Verde"
$ws.Cells.Item(141, 2).Value = 328
$ws.Cells.Item(141, 4).Value = 84
$ws.Cells.Item(141, 5).Value = 241
$ws.Cells.Item(141, 8).Value = 3

# Row 142 -> now "Republica de Africa Central"
$ws.Cells.Item(142, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(142, 2).Value = 327
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 13
$ws.Cells.Item(142, 5).Value = 314

# Row 143 -> now "Vietnam"
$ws.Cells.Item(143, 1).Value = "Vietnam"
$ws.Cells.Item(143, 3).Value = 4
$ws.Cells.Item(143, 4).Value = 263
$ws.Cells.Item(143, 5).Value = 61
$ws.Cells.Item(143, 8).Value = 0

# Row 144 -> now "Montenegro"
$ws.Cells.Item(144, 1).Value = "Montenegro"
$ws.Cells.Item(144, 2).Value = 324
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 311
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 8).Value = 9

# Row 145 -> now "Madagascar"
$ws.Cells.Item(145, 1).Value = "Madagascar"
$ws.Cells.Item(145, 2).Value = 322
$ws.Cells.Item(145, 3).Value = 18
$ws.Cells.Item(145, 4).Value = 119
$ws.Cells.Item(145, 5).Value = 202
$ws.Cells.Item(145, 8).Value = 1

# --- Yemen (row 161) ---------------------------------------------------------
$ws.Cells.Item(161, 2).Value = 130
$ws.Cells.Item(161, 3).Value = 2
$ws.Cells.Item(161, 5).Value = 109
